$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

# Remove rows 3..12, keep header (A1) and A2
$ws.Range("A3:A12").EntireRow.Delete() | Out-Null

# Set A2 value and hyperlink
$ws.Range("A2").Value = "https://nursesopenings.com/"
$ws.Hyperlinks.Add($ws.Range("A2"), "https://nursesopenings.com/") | Out-Null

$ws.Range("A2:A12").Select() | Out-Null
